$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 271.8421
$ws.Range("I33").Value = 261.94446
$ws.Range("J33").Value = 450
$ws.Range("K33").Value = 261.94446
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = -32.94445999999999
$ws.Range("N33").Value = -908

$ws.Range("H37").Value = 5000
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15252

$ws.Range("H62").Value = 1709.5
$ws.Range("J62").Value = 1438
$ws.Range("L62").Value = 1438
$ws.Range("N62").Value = -2686

$ws.Range("H65").Value = 1709.5
$ws.Range("J65").Value = 1438
$ws.Range("L65").Value = 7190
$ws.Range("N65").Value = -13430

$ws.Range("H76").Value = 3499.95
$ws.Range("I76").Value = 3294.0588
$ws.Range("J76").Value = 4666.6665
$ws.Range("K76").Value = 3294.0588
$ws.Range("L76").Value = 4666.6665
$ws.Range("M76").Value = -2979.0588
$ws.Range("N76").Value = -5296.6665

$ws.Range("H79").Value = 3499.95
$ws.Range("I79").Value = 3294.0588
$ws.Range("J79").Value = 4666.6665
$ws.Range("K79").Value = 3294.0588
$ws.Range("L79").Value = 4666.6665
$ws.Range("M79").Value = -2202.0588
$ws.Range("N79").Value = -6850.6665

$ws.Range("H81").Value = 38328
$ws.Range("J81").Value = 38328
$ws.Range("L81").Value = 38328
$ws.Range("N81").Value = -40324

$ws.Range("H84").Value = 38328
$ws.Range("J84").Value = 38328
$ws.Range("L84").Value = 114984
$ws.Range("N84").Value = -124968

$ws.Range("H88").Value = 1934.3334
$ws.Range("I88").Value = 2003
$ws.Range("J88").Value = 1900
$ws.Range("K88").Value = 2003
$ws.Range("L88").Value = 1900
$ws.Range("M88").Value = -1597
$ws.Range("N88").Value = -2712

$ws.Range("H91").Value = 1934.3334
$ws.Range("I91").Value = 2003
$ws.Range("J91").Value = 1900
$ws.Range("K91").Value = 2003
$ws.Range("L91").Value = 1900
$ws.Range("M91").Value = -599
$ws.Range("N91").Value = -4708

$ws.Range("H98").Value = 1010
$ws.Range("I98").Value = 1103.5714
$ws.Range("J98").Value = 791.6667
$ws.Range("K98").Value = 1103.5714
$ws.Range("L98").Value = 791.6667
$ws.Range("M98").Value = 394.4286
$ws.Range("N98").Value = -3787.6667

$ws.Range("H105").Value = 39671
$ws.Range("J105").Value = 39671
$ws.Range("L105").Value = 39671
$ws.Range("N105").Value = -46659

$ws.Range("H106").Value = 1887.1428
$ws.Range("I106").Value = 1887.1428
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1887.1428
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1256.1428
$ws.Range("N106").ClearContents()

$ws.Range("H108").Value = 41172
$ws.Range("J108").Value = 41172
$ws.Range("L108").Value = 41172
$ws.Range("N108").Value = -48852

$ws.Range("H113").Value = 2947.389
$ws.Range("I113").Value = 2425
$ws.Range("J113").Value = 3096.6428
$ws.Range("K113").Value = 2425
$ws.Range("L113").Value = 3096.6428
$ws.Range("M113").Value = 829
$ws.Range("N113").Value = -9604.6428

$ws.Range("H117").Value = 51759.8
$ws.Range("J117").Value = 51759.8
$ws.Range("L117").Value = 51759.8
$ws.Range("N117").Value = -60937.8

$ws.Range("H120").Value = 77893.336
$ws.Range("J120").Value = 77893.336
$ws.Range("L120").Value = 77893.336
$ws.Range("N120").Value = -87569.336

$ws.Range("H122").Value = 1010
$ws.Range("I122").Value = 1103.5714
$ws.Range("J122").Value = 791.6667
$ws.Range("K122").Value = 3310.7142
$ws.Range("L122").Value = 2375.0001
$ws.Range("M122").Value = -860.7142000000003
$ws.Range("N122").Value = -7275.0001

$ws.Range("H123").Value = 35266.332
$ws.Range("J123").Value = 35266.332
$ws.Range("L123").Value = 35266.332
$ws.Range("N123").Value = -45066.332

$ws.Range("H124").Value = 78900
$ws.Range("J124").Value = 78900
$ws.Range("L124").Value = 78900
$ws.Range("N124").Value = -88720

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1828.05
$ws.Range("I45").Value = 1784.3125
$ws.Range("K45").Value = 1784.3125
$ws.Range("M45").Value = -1407.3125

$ws.Range("H63").Value = 4179.3
$ws.Range("I63").Value = 3967.1667
$ws.Range("J63").Value = 4497.5
$ws.Range("K63").Value = 3967.1667
$ws.Range("L63").Value = 4497.5
$ws.Range("M63").Value = -3281.1667
$ws.Range("N63").Value = -5869.5

$ws.Range("H66").Value = 4179.3
$ws.Range("I66").Value = 3967.1667
$ws.Range("J66").Value = 4497.5
$ws.Range("K66").Value = 19835.8335
$ws.Range("L66").Value = 22487.5
$ws.Range("M66").Value = -16403.8335
$ws.Range("N66").Value = -29351.5

$ws.Range("H118").Value = 29983.166
$ws.Range("J118").Value = 29983.166
$ws.Range("L118").Value = 29983.166
$ws.Range("N118").Value = -33297.166

$ws.Range("H122").Value = 4033294.8
$ws.Range("I122").Value = 1034.591
$ws.Range("J122").Value = 13889931
$ws.Range("K122").Value = 3103.773
$ws.Range("L122").Value = 41669793
$ws.Range("M122").Value = -653.7729999999997
$ws.Range("N122").Value = -41674693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 37000
$ws.Range("J51").Value = 37000
$ws.Range("L51").Value = 37000
$ws.Range("N51").Value = -37982

$ws.Range("H94").Value = 1585.7142
$ws.Range("I94").Value = 1350
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1350
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -899
$ws.Range("N94").Value = -3902

$ws.Range("H99").Value = 1567.619
$ws.Range("I99").Value = 1601
$ws.Range("J99").Value = 1537.2727
$ws.Range("K99").Value = 1601
$ws.Range("L99").Value = 1537.2727
$ws.Range("M99").Value = -103
$ws.Range("N99").Value = -4533.2727

$ws.Range("H123").Value = 75382
$ws.Range("J123").Value = 75382
$ws.Range("L123").Value = 75382
$ws.Range("N123").Value = -85182

$ws.Range("H134").Value = 59817.055
$ws.Range("I134").Value = 4741.0835
$ws.Range("J134").Value = 169969
$ws.Range("K134").Value = 14223.2505
$ws.Range("L134").Value = 509907
$ws.Range("M134").Value = -11688.2505
$ws.Range("N134").Value = -514977

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7495.0415
$ws.Range("I31").Value = 10664.25
$ws.Range("K31").Value = 10664.25
$ws.Range("M31").Value = -10369.25

$ws.Range("H34").Value = 7495.0415
$ws.Range("I34").Value = 10664.25
$ws.Range("K34").Value = 10664.25
$ws.Range("M34").Value = -10462.25

$ws.Range("H122").Value = 8609
$ws.Range("I122").Value = 3232.2354
$ws.Range("J122").Value = 100014
$ws.Range("K122").Value = 9696.706200000001
$ws.Range("L122").Value = 300042
$ws.Range("M122").Value = -7246.706200000001
$ws.Range("N122").Value = -304942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 10027.75
$ws.Range("J97").Value = 14397.875
$ws.Range("L97").Value = 43193.625
$ws.Range("N97").Value = -44185.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10429.429
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 3834.3333
$ws.Range("K80").Value = 50000
$ws.Range("L80").Value = 3834.3333
$ws.Range("M80").Value = -49002
$ws.Range("N80").Value = -5830.3333

$ws.Range("H83").Value = 10429.429
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 3834.3333
$ws.Range("K83").Value = 250000
$ws.Range("L83").Value = 19171.6665
$ws.Range("M83").Value = -245008
$ws.Range("N83").Value = -29155.6665

$ws.Range("H122").Value = 3923.1667
$ws.Range("I122").Value = 4045.111
$ws.Range("J122").Value = 3557.3333
$ws.Range("K122").Value = 12135.333
$ws.Range("L122").Value = 10671.9999
$ws.Range("M122").Value = -9685.332999999999
$ws.Range("N122").Value = -15571.9999

$ws.Range("H126").Value = 2825.3704
$ws.Range("I126").Value = 1887.5
$ws.Range("J126").Value = 3220.2632
$ws.Range("K126").Value = 5662.5
$ws.Range("L126").Value = 9660.7896
$ws.Range("M126").Value = -3192.5
$ws.Range("N126").Value = -14600.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1117.8572
$ws.Range("I22").Value = 1064.6
$ws.Range("J22").Value = 1251
$ws.Range("K22").Value = 1064.6
$ws.Range("L22").Value = 1251
$ws.Range("M22").Value = -769.5999999999999
$ws.Range("N22").Value = -1841

$ws.Range("H27").Value = 1117.8572
$ws.Range("I27").Value = 1064.6
$ws.Range("J27").Value = 1251
$ws.Range("K27").Value = 1064.6
$ws.Range("L27").Value = 1251
$ws.Range("M27").Value = -957.5999999999999
$ws.Range("N27").Value = -1465

$ws.Range("H40").Value = 4797.2666
$ws.Range("I40").Value = 4612.231
$ws.Range("K40").Value = 4612.231
$ws.Range("M40").Value = -4476.231

$ws.Range("H46").Value = 769.2308
$ws.Range("I46").Value = 312.5
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 312.5
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -124.5
$ws.Range("N46").Value = -1876

$ws.Range("H93").Value = 467.875
$ws.Range("I93").Value = 475.92307
$ws.Range("J93").Value = 458.36365
$ws.Range("K93").Value = 475.92307
$ws.Range("L93").Value = 458.36365
$ws.Range("M93").Value = 772.0769299999999
$ws.Range("N93").Value = -2954.36365

$ws.Range("H122").Value = 6733.1724
$ws.Range("I122").Value = 6353.8667
$ws.Range("J122").Value = 8046.154
$ws.Range("K122").Value = 19061.6001
$ws.Range("L122").Value = 24138.462
$ws.Range("M122").Value = -16611.6001
$ws.Range("N122").Value = -29038.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 59266.668
$ws.Range("J109").Value = 59266.668
$ws.Range("L109").Value = 59266.668
$ws.Range("N109").Value = -62040.668
